# Generate Report for Handback
# - marks the zh-cn / de-de rows as handed back (status text + handback file +
#   handback datetime), adds "Latest Target File" hyperlinks, and widens a
#   couple of columns so the new (longer) text fits.

$wb = $excel.ActiveWorkbook

$AMD_URL = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb1826a5a779122265f7a89a1ca83b73d614b95f/e2e/a.md"
$HANDED_BACK = "Handed back: in sync with en-US"
$HYPERLINK_COLOR = 15570276   # OLE/BGR encoding of RGB(0x64,0x95,0xED) == FF6495ED

function Set-HandoffHyperlink($ws, $cellRef) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $AMD_URL, $null, $null, "a.md") | Out-Null
    $r = $ws.Range($cellRef)
    $r.Font.Underline = $true
    $r.Font.Color = $HYPERLINK_COLOR
    $r.Font.Name = "Calibri"
}

# ---------------------------------------------------------------------------
# Overview sheet: the status text itself lives in the shared string that
# E2/F2/E3/F3 already point to, so just rewriting their values in place
# propagates to every cell that shares that string (same thing happens on the
# zh-cn / de-de sheets' Status column below).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $HANDED_BACK
$wsOverview.Range("F2").Value = $HANDED_BACK
$wsOverview.Range("E3").Value = $HANDED_BACK
$wsOverview.Range("F3").Value = $HANDED_BACK

$wsOverview.Columns.Item(5).ColumnWidth = 30
$wsOverview.Columns.Item(6).ColumnWidth = 30

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $HANDED_BACK
$wsZh.Range("C3").Value = $HANDED_BACK

Set-HandoffHyperlink $wsZh "I2"
Set-HandoffHyperlink $wsZh "I3"

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-28 08:36:24"
$wsZh.Range("K3").Value = "2016-08-28 08:36:24"

$wsZh.Columns.Item(3).ColumnWidth = 30
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $HANDED_BACK
$wsDe.Range("C3").Value = $HANDED_BACK

Set-HandoffHyperlink $wsDe "I2"
Set-HandoffHyperlink $wsDe "I3"

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-28 08:36:30"
$wsDe.Range("K3").Value = "2016-08-28 08:36:30"

$wsDe.Columns.Item(3).ColumnWidth = 30
$wsDe.Columns.Item(10).ColumnWidth = 40
